$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells keep their text format (values contain dotted numbers / padded percents
# that must not be reinterpreted as numeric/date types).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.301.47"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +2.86%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.812.09"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +3.85%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.14%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4358"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +2.39%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3665"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +1.47%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "44.87"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.07%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07682"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +2.88%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.142"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +2.32%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9997"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.34%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.00"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +2.33%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.304"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +3.29%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.525"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +4.41%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.827.73"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +5.17%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "95.36"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +8.86%  "
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.22%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06535"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +4.84%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9998"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.03%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.43"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +2.99%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.247"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +2.10%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.316.46"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +2.76%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.36%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -10.44%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "161.94"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +6.89%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.74"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.43%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.024.20"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +4.57%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.279"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -3.65%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "128.92"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.76%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.209"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.10%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.962"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +4.66%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09187"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.56%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.495"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -5.01%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "13.02"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +2.70%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02348"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +2.16%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.191"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +2.08%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +1.42%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.6588"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +2.38%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +2.03%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.07%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.139"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +2.63%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.432"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.22%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9990"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.06%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.65%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6114"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +3.64%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.742"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.50%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "125.88"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.14%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.018"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +2.93%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.156"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +3.21%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06995"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +1.85%  "
